$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 updates
$ws.Range("G11").Value = 2.1
$ws.Range("I11").Value = 3.8
$ws.Range("Q11").Value = 2.15
$ws.Range("R11").Value = 1.67
$ws.Range("U11").Value = 1.83
$ws.Range("V11").Value = 1.83
$ws.Range("W11").Value = 7
$ws.Range("X11").Value = 9.5
$ws.Range("Z11").Value = 19
$ws.Range("AB11").Value = 29
$ws.Range("AI11").Value = 13
$ws.Range("AM11").Value = 301
$ws.Range("AO11").Value = 12
$ws.Range("AV11").Value = 51
$ws.Range("AY11").Value = 29

# Row 12 updates
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 9
